$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (sheet1) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 12868
$ws1.Range("F3").Value = 632
$ws1.Range("F6").Value = 324
$ws1.Range("F7").Value = 406
$ws1.Range("F8").Value = 236
$ws1.Range("F9").Value = 12915
$ws1.Range("F11").Value = 26
$ws1.Range("F12").Value = 5256
$ws1.Range("F13").Value = 548
$ws1.Range("F14").Value = 20
$ws1.Range("F15").Value = 14
$ws1.Range("F16").Value = 31
$ws1.Range("F20").Value = 682
$ws1.Range("F21").Value = 2858
$ws1.Range("F22").Value = 6191
$ws1.Range("F23").Value = 1161
$ws1.Range("F24").Value = 3629
$ws1.Range("F25").Value = 222

# --- Sheet "全部类型" (sheet4) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 12868
$ws4.Range("F3").Value = 632
$ws4.Range("F6").Value = 324
$ws4.Range("F8").Value = 406
$ws4.Range("F9").Value = 236
$ws4.Range("F10").Value = 0
$ws4.Range("F12").Value = 26
$ws4.Range("F13").Value = 5256
$ws4.Range("F14").Value = 548
$ws4.Range("F15").Value = 20
$ws4.Range("F16").Value = 14
$ws4.Range("F17").Value = 31
$ws4.Range("F21").Value = 682
$ws4.Range("F22").Value = 2858
$ws4.Range("F24").Value = 6191
$ws4.Range("F25").Value = 1161
$ws4.Range("F26").Value = 3629
$ws4.Range("F27").Value = 222

$wb.Save()
